# ACTUALIZACION DE AVANCES DE LA SEMANA
# Updates the "Enero" sheet with the latest weekly progress scores
# (column D, I, N, S and X were still showing 0 for "Day 3" of each
# week; also the first-week Monday score in column B moved 7 -> 8).
# All of the AVERAGE() formulas down-stream (F:F6, K:K6, P:P6, U:U6,
# Z:Z6, AA:AF and the totals row 7) recalculate automatically, as do
# the cached chart values on the sheet's dashboard charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enero")

# Semana 1 (row 4)
$ws.Range("B4").Value = 8
$ws.Range("D4").Value = 7
$ws.Range("I4").Value = 7
$ws.Range("N4").Value = 8
$ws.Range("S4").Value = 9
$ws.Range("X4").Value = 8

# Semana 2 (row 5)
$ws.Range("D5").Value = 7
$ws.Range("I5").Value = 7
$ws.Range("N5").Value = 7
$ws.Range("S5").Value = 9
$ws.Range("X5").Value = 7

# Semana 3 (row 6)
$ws.Range("D6").Value = 7
$ws.Range("I6").Value = 7
$ws.Range("N6").Value = 7
$ws.Range("S6").Value = 9
$ws.Range("X6").Value = 8

# Restore the cursor to where the author left off editing.
$ws.Activate()
$ws.Range("AL11").Select()

$wb.Save()
